$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 2864.8462
$ws.Range("J5").Value2 = 5141
$ws.Range("L5").Value2 = 5141
$ws.Range("N5").Value2 = -5371
$ws.Range("H18").Value2 = 1615
$ws.Range("I18").Value2 = 421.5
$ws.Range("J18").Value2 = 4002
$ws.Range("K18").Value2 = 421.5
$ws.Range("L18").Value2 = 4002
$ws.Range("M18").Value2 = -137.5
$ws.Range("N18").Value2 = -4570
$ws.Range("H19").Value2 = 1187.8
$ws.Range("I19").Value2 = 758.26666
$ws.Range("J19").Value2 = 1617.3334
$ws.Range("K19").Value2 = 758.26666
$ws.Range("L19").Value2 = 1617.3334
$ws.Range("M19").Value2 = -583.26666
$ws.Range("N19").Value2 = -1967.3334
$ws.Range("H32").Value2 = 11303.467
$ws.Range("I32").Value2 = 8173.5
$ws.Range("J32").Value2 = 12441.637
$ws.Range("K32").Value2 = 8173.5
$ws.Range("L32").Value2 = 12441.637
$ws.Range("M32").Value2 = -7847.5
$ws.Range("N32").Value2 = -13093.637
$ws.Range("H40").Value2 = 2999.8572
$ws.Range("I40").Value2 = 2399.4
$ws.Range("J40").Value2 = 4501
$ws.Range("K40").Value2 = 2399.4
$ws.Range("L40").Value2 = 4501
$ws.Range("M40").Value2 = -2224.4
$ws.Range("N40").Value2 = -4851
$ws.Range("H51").Value2 = 15956.714
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 15956.714
$ws.Range("K51").Value2 = 0
$ws.Range("L51").ClearContents()
$ws.Range("M51").Value2 = 15956.714
$ws.Range("N51").Value2 = -16924.714
$ws.Range("H58").Value2 = 400
$ws.Range("J58").Value2 = 0
$ws.Range("L58").Value2 = 0
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value2 = 5890442.5
$ws.Range("I62").Value2 = 21853934
$ws.Range("J62").Value2 = 9156
$ws.Range("K62").Value2 = 21853934
$ws.Range("L62").Value2 = 9156
$ws.Range("M62").Value2 = -21853310
$ws.Range("N62").Value2 = -10404
$ws.Range("H64").Value2 = 23824278
$ws.Range("I64").Value2 = 55565956
$ws.Range("K64").Value2 = 55565956
$ws.Range("M64").Value2 = -55565708
$ws.Range("H65").Value2 = 5890442.5
$ws.Range("I65").Value2 = 21853934
$ws.Range("J65").Value2 = 9156
$ws.Range("K65").Value2 = 109269670
$ws.Range("L65").Value2 = 45780
$ws.Range("M65").Value2 = -109266550
$ws.Range("N65").Value2 = -52020
$ws.Range("H67").Value2 = 23824278
$ws.Range("I67").Value2 = 55565956
$ws.Range("K67").Value2 = 55565956
$ws.Range("M67").Value2 = -55565098
$ws.Range("H70").Value2 = 2993
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 2993
$ws.Range("K70").Value2 = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value2 = 8979
$ws.Range("N70").Value2 = -9519
$ws.Range("H73").Value2 = 2993
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 2993
$ws.Range("K73").Value2 = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value2 = 8979
$ws.Range("N73").Value2 = -10851
$ws.Range("H76").Value2 = 101007736
$ws.Range("I76").Value2 = 144290340
$ws.Range("K76").Value2 = 144290340
$ws.Range("M76").Value2 = -144290025
$ws.Range("H79").Value2 = 101007736
$ws.Range("I79").Value2 = 144290340
$ws.Range("K79").Value2 = 144290340
$ws.Range("M79").Value2 = -144289248
$ws.Range("H86").Value2 = 2001267.2
$ws.Range("I86").Value2 = 2001267.2
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 2001267.2
$ws.Range("L86").Value2 = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value2 = -2000144.2
$ws.Range("H87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value2 = 0
$ws.Range("H88").Value2 = 10000
$ws.Range("J88").Value2 = 10000
$ws.Range("L88").Value2 = 10000
$ws.Range("N88").Value2 = -10812
$ws.Range("H89").Value2 = 2001267.2
$ws.Range("I89").Value2 = 2001267.2
$ws.Range("J89").Value2 = 0
$ws.Range("K89").Value2 = 10006336
$ws.Range("L89").Value2 = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value2 = -10000720
$ws.Range("H90").Value2 = 0
$ws.Range("J90").Value2 = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value2 = 0
$ws.Range("H91").Value2 = 10000
$ws.Range("J91").Value2 = 10000
$ws.Range("L91").Value2 = 10000
$ws.Range("N91").Value2 = -12808
$ws.Range("H103").Value2 = 511.44446
$ws.Range("I103").Value2 = 505
$ws.Range("J103").Value2 = 516.6
$ws.Range("K103").Value2 = 1515
$ws.Range("L103").Value2 = 1549.8
$ws.Range("M103").Value2 = -929
$ws.Range("N103").Value2 = -2721.8
$ws.Range("H116").Value2 = 5855192.5
$ws.Range("I116").Value2 = 9265528
$ws.Range("K116").Value2 = 9265528
$ws.Range("M116").Value2 = -9262086
$ws.Range("H132").Value2 = 384569.25
$ws.Range("I132").Value2 = 463798.34
$ws.Range("K132").Value2 = 1391395.02
$ws.Range("M132").Value2 = -1388865.02
$ws.Range("H137").Value2 = 7054.6665
$ws.Range("I137").Value2 = 5288.5625
$ws.Range("J137").Value2 = 10586.875
$ws.Range("K137").Value2 = 15865.6875
$ws.Range("L137").Value2 = 31760.625
$ws.Range("M137").Value2 = -13315.6875
$ws.Range("N137").Value2 = -36860.625
$ws.Range("H138").Value2 = 3339.5894
$ws.Range("I138").Value2 = 1735.3684
$ws.Range("J138").Value2 = 4163.3784
$ws.Range("K138").Value2 = 5206.1052
$ws.Range("L138").Value2 = 12490.1352
$ws.Range("M138").Value2 = -66.10519999999997
$ws.Range("N138").Value2 = -22770.1352
$ws.Range("H141").Value2 = 2999.5
$ws.Range("J141").Value2 = 4999
$ws.Range("L141").Value2 = 14997
$ws.Range("N141").Value2 = -25357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 127.55556
$ws.Range("I4").Value2 = 127.55556
$ws.Range("K4").Value2 = 127.55556
$ws.Range("M4").Value2 = -11.55556
$ws.Range("H32").Value2 = 1985468.6
$ws.Range("I32").Value2 = 996.83673
$ws.Range("K32").Value2 = 996.83673
$ws.Range("M32").Value2 = -709.83673
$ws.Range("H45").Value2 = 2358.1333
$ws.Range("I45").Value2 = 2170.818
$ws.Range("J45").Value2 = 2873.25
$ws.Range("K45").Value2 = 2170.818
$ws.Range("L45").Value2 = 2873.25
$ws.Range("M45").Value2 = -1793.818
$ws.Range("N45").Value2 = -3627.25
$ws.Range("H61").Value2 = 5475.421
$ws.Range("I61").Value2 = 7113.1177
$ws.Range("K61").Value2 = 7113.1177
$ws.Range("M61").Value2 = -6901.1177
$ws.Range("H74").Value2 = 2793.139
$ws.Range("I74").Value2 = 2691.4614
$ws.Range("J74").Value2 = 2850.6086
$ws.Range("K74").Value2 = 2691.4614
$ws.Range("L74").Value2 = 2850.6086
$ws.Range("M74").Value2 = -1817.4614
$ws.Range("N74").Value2 = -4598.6086
$ws.Range("H77").Value2 = 2793.139
$ws.Range("I77").Value2 = 2691.4614
$ws.Range("J77").Value2 = 2850.6086
$ws.Range("K77").Value2 = 13457.307
$ws.Range("L77").Value2 = 14253.043
$ws.Range("M77").Value2 = -9089.307000000001
$ws.Range("N77").Value2 = -22989.043
$ws.Range("H97").Value2 = 381.57895
$ws.Range("I97").Value2 = 306.06668
$ws.Range("J97").Value2 = 664.75
$ws.Range("K97").Value2 = 306.06668
$ws.Range("L97").Value2 = 664.75
$ws.Range("M97").Value2 = 189.93332
$ws.Range("N97").Value2 = -1656.75
$ws.Range("H122").Value2 = 2755.8235
$ws.Range("I122").Value2 = 1735.25
$ws.Range("K122").Value2 = 5205.75
$ws.Range("M122").Value2 = -2755.75
$ws.Range("H136").Value2 = 5475.421
$ws.Range("I136").Value2 = 7113.1177
$ws.Range("K136").Value2 = 21339.3531
$ws.Range("M136").Value2 = -18789.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 6494846
$ws.Range("I20").Value2 = 12988204
$ws.Range("J20").Value2 = 1488.2727
$ws.Range("K20").Value2 = 12988204
$ws.Range("L20").Value2 = 1488.2727
$ws.Range("M20").Value2 = -12987957
$ws.Range("N20").Value2 = -1982.2727
$ws.Range("H22").Value2 = 1894.5
$ws.Range("I22").Value2 = 1894.5
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 1894.5
$ws.Range("L22").Value2 = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value2 = -1721.5
$ws.Range("H86").Value2 = 5542.636
$ws.Range("I86").Value2 = 2856.6667
$ws.Range("J86").Value2 = 11298.286
$ws.Range("K86").Value2 = 2856.6667
$ws.Range("L86").Value2 = 11298.286
$ws.Range("M86").Value2 = -1733.6667
$ws.Range("N86").Value2 = -13544.286
$ws.Range("H89").Value2 = 5542.636
$ws.Range("I89").Value2 = 2856.6667
$ws.Range("J89").Value2 = 11298.286
$ws.Range("K89").Value2 = 14283.3335
$ws.Range("L89").Value2 = 56491.43
$ws.Range("M89").Value2 = -8667.333500000001
$ws.Range("N89").Value2 = -67723.42999999999
$ws.Range("H94").Value2 = 5059.5
$ws.Range("I94").Value2 = 641.1539
$ws.Range("K94").Value2 = 641.1539
$ws.Range("M94").Value2 = -190.1539
$ws.Range("H97").Value2 = 10462.833
$ws.Range("I97").Value2 = 10462.833
$ws.Range("K97").Value2 = 10462.833
$ws.Range("M97").Value2 = -9471.833000000001
$ws.Range("H105").Value2 = 6913.125
$ws.Range("I105").Value2 = 3404.75
$ws.Range("K105").Value2 = 3404.75
$ws.Range("M105").Value2 = -1657.75
$ws.Range("H134").Value2 = 789274.25
$ws.Range("I134").Value2 = 997448.6
$ws.Range("J134").Value2 = 20630.54
$ws.Range("K134").Value2 = 2992345.8
$ws.Range("L134").Value2 = 61891.62
$ws.Range("M134").Value2 = -2989810.8
$ws.Range("N134").Value2 = -66961.62

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 3450
$ws.Range("I4").Value2 = 5000
$ws.Range("J4").Value2 = 2933.3333
$ws.Range("K4").Value2 = 5000
$ws.Range("L4").Value2 = 2933.3333
$ws.Range("M4").Value2 = -4888
$ws.Range("N4").Value2 = -3157.3333
$ws.Range("H22").Value2 = 812471.5600000001
$ws.Range("I22").Value2 = 1116717.1
$ws.Range("K22").Value2 = 1116717.1
$ws.Range("M22").Value2 = -1116367.1
$ws.Range("H28").Value2 = 13189.857
$ws.Range("J28").Value2 = 13189.857
$ws.Range("L28").Value2 = 13189.857
$ws.Range("N28").Value2 = -13679.857
$ws.Range("H31").Value2 = 3614.4119
$ws.Range("I31").Value2 = 1060
$ws.Range("J31").Value2 = 5552.241
$ws.Range("K31").Value2 = 1060
$ws.Range("L31").Value2 = 5552.241
$ws.Range("M31").Value2 = -765
$ws.Range("N31").Value2 = -6142.241
$ws.Range("H34").Value2 = 3614.4119
$ws.Range("I34").Value2 = 1060
$ws.Range("J34").Value2 = 5552.241
$ws.Range("K34").Value2 = 1060
$ws.Range("L34").Value2 = 5552.241
$ws.Range("M34").Value2 = -858
$ws.Range("N34").Value2 = -5956.241
$ws.Range("H58").Value2 = 47628236
$ws.Range("I58").Value2 = 76930250
$ws.Range("K58").Value2 = 76930250
$ws.Range("M58").Value2 = -76930047
$ws.Range("H59").Value2 = 10000
$ws.Range("I59").Value2 = 10000
$ws.Range("K59").Value2 = 10000
$ws.Range("M59").Value2 = -8855
$ws.Range("H69").Value2 = 10545.25
$ws.Range("I69").Value2 = 7393.6665
$ws.Range("K69").Value2 = 7393.6665
$ws.Range("M69").Value2 = -6644.6665
$ws.Range("H72").Value2 = 10545.25
$ws.Range("I72").Value2 = 7393.6665
$ws.Range("K72").Value2 = 22180.9995
$ws.Range("M72").Value2 = -18436.9995
$ws.Range("H86").Value2 = 10087.5
$ws.Range("I86").Value2 = 9646.777
$ws.Range("J86").Value2 = 10392.615
$ws.Range("K86").Value2 = 9646.777
$ws.Range("L86").Value2 = 10392.615
$ws.Range("M86").Value2 = -8523.777
$ws.Range("N86").Value2 = -12638.615
$ws.Range("H89").Value2 = 10087.5
$ws.Range("I89").Value2 = 9646.777
$ws.Range("J89").Value2 = 10392.615
$ws.Range("K89").Value2 = 48233.885
$ws.Range("L89").Value2 = 51963.075
$ws.Range("M89").Value2 = -42617.885
$ws.Range("N89").Value2 = -63195.075
$ws.Range("H99").Value2 = 15156684
$ws.Range("J99").Value2 = 4998
$ws.Range("L99").Value2 = 4998
$ws.Range("N99").Value2 = -7994
$ws.Range("H105").Value2 = 35717030
$ws.Range("I105").Value2 = 43480204
$ws.Range("K105").Value2 = 43480204
$ws.Range("M105").Value2 = -43478457
$ws.Range("H126").Value2 = 15156684
$ws.Range("J126").Value2 = 4998
$ws.Range("L126").Value2 = 14994
$ws.Range("N126").Value2 = -19934
$ws.Range("H132").Value2 = 24047.805
$ws.Range("I132").Value2 = 26053.875
$ws.Range("K132").Value2 = 78161.625
$ws.Range("M132").Value2 = -75631.625
$ws.Range("H134").Value2 = 41678356
$ws.Range("I134").Value2 = 66677868
$ws.Range("K134").Value2 = 200033604
$ws.Range("M134").Value2 = -200031069
$ws.Range("H136").Value2 = 47628236
$ws.Range("I136").Value2 = 76930250
$ws.Range("K136").Value2 = 230790750
$ws.Range("M136").Value2 = -230788200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 122.333336
$ws.Range("I2").Value2 = 120.5
$ws.Range("J2").Value2 = 126
$ws.Range("K2").Value2 = 723
$ws.Range("L2").Value2 = 756
$ws.Range("M2").Value2 = -610
$ws.Range("N2").Value2 = -982
$ws.Range("H9").Value2 = 145483.33
$ws.Range("I9").Value2 = 235750
$ws.Range("J9").Value2 = 100350
$ws.Range("K9").Value2 = 707250
$ws.Range("L9").Value2 = 301050
$ws.Range("M9").Value2 = -707026
$ws.Range("N9").Value2 = -301498
$ws.Range("H14").Value2 = 569.375
$ws.Range("I14").Value2 = 569.375
$ws.Range("K14").Value2 = 1708.125
$ws.Range("M14").Value2 = -1535.125
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("K20").Value2 = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value2 = 0
$ws.Range("H21").Value2 = 2384.5
$ws.Range("I21").Value2 = 770
$ws.Range("J21").Value2 = 3999
$ws.Range("K21").Value2 = 2310
$ws.Range("L21").Value2 = 11997
$ws.Range("M21").Value2 = -2137
$ws.Range("N21").Value2 = -12343
$ws.Range("H22").Value2 = 1452.25
$ws.Range("I22").Value2 = 255
$ws.Range("J22").Value2 = 2649.5
$ws.Range("K22").Value2 = 765
$ws.Range("L22").Value2 = 7948.5
$ws.Range("M22").Value2 = -596
$ws.Range("N22").Value2 = -8286.5
$ws.Range("H27").Value2 = 1452.25
$ws.Range("I27").Value2 = 255
$ws.Range("J27").Value2 = 2649.5
$ws.Range("K27").Value2 = 765
$ws.Range("L27").Value2 = 7948.5
$ws.Range("M27").Value2 = -663
$ws.Range("N27").Value2 = -8152.5
$ws.Range("H29").Value2 = 240.88889
$ws.Range("I29").Value2 = 182.8
$ws.Range("J29").Value2 = 313.5
$ws.Range("K29").Value2 = 548.4000000000001
$ws.Range("L29").Value2 = 940.5
$ws.Range("M29").Value2 = -271.4000000000001
$ws.Range("N29").Value2 = -1494.5
$ws.Range("H33").Value2 = 46.846153
$ws.Range("I33").Value2 = 31.583334
$ws.Range("K33").Value2 = 189.500004
$ws.Range("M33").Value2 = 93.49999600000001
$ws.Range("H35").Value2 = 795.4286
$ws.Range("I35").Value2 = 1300
$ws.Range("J35").Value2 = 711.3333
$ws.Range("K35").Value2 = 3900
$ws.Range("L35").Value2 = 2133.9999
$ws.Range("M35").Value2 = -3612
$ws.Range("N35").Value2 = -2709.9999
$ws.Range("H38").Value2 = 194.10715
$ws.Range("J38").Value2 = 599.125
$ws.Range("L38").Value2 = 1797.375
$ws.Range("N38").Value2 = -2491.375
$ws.Range("H64").Value2 = 14649.667
$ws.Range("I64").Value2 = 14500
$ws.Range("J64").Value2 = 14724.5
$ws.Range("K64").Value2 = 43500
$ws.Range("L64").Value2 = 44173.5
$ws.Range("M64").Value2 = -43230
$ws.Range("N64").Value2 = -44713.5
$ws.Range("H67").Value2 = 14649.667
$ws.Range("I67").Value2 = 14500
$ws.Range("J67").Value2 = 14724.5
$ws.Range("K67").Value2 = 43500
$ws.Range("L67").Value2 = 44173.5
$ws.Range("M67").Value2 = -42564
$ws.Range("N67").Value2 = -46045.5
$ws.Range("H68").Value2 = 120804.64
$ws.Range("J68").Value2 = 130002.38
$ws.Range("L68").Value2 = 390007.14
$ws.Range("N68").Value2 = -391629.14
$ws.Range("H69").Value2 = 25000
$ws.Range("I69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("M69").ClearContents()
$ws.Range("H71").Value2 = 120804.64
$ws.Range("J71").Value2 = 130002.38
$ws.Range("L71").Value2 = 1170021.42
$ws.Range("N71").Value2 = -1178133.42
$ws.Range("H72").Value2 = 25000
$ws.Range("I72").Value2 = 0
$ws.Range("K72").Value2 = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value2 = 0
$ws.Range("I80").Value2 = 0
$ws.Range("K80").Value2 = 0
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value2 = 2179.4
$ws.Range("I81").Value2 = 1632.6666
$ws.Range("J81").Value2 = 2999.5
$ws.Range("K81").Value2 = 4897.9998
$ws.Range("L81").Value2 = 8998.5
$ws.Range("M81").Value2 = -3774.9998
$ws.Range("N81").Value2 = -11244.5
$ws.Range("H82").Value2 = 22916.666
$ws.Range("I82").Value2 = 0
$ws.Range("J82").Value2 = 22916.666
$ws.Range("K82").Value2 = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").Value2 = 68749.99800000001
$ws.Range("N82").Value2 = -69561.99800000001
$ws.Range("H83").Value2 = 0
$ws.Range("I83").Value2 = 0
$ws.Range("K83").Value2 = 0
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value2 = 2179.4
$ws.Range("I84").Value2 = 1632.6666
$ws.Range("J84").Value2 = 2999.5
$ws.Range("K84").Value2 = 14693.9994
$ws.Range("L84").Value2 = 26995.5
$ws.Range("M84").Value2 = -9077.999400000001
$ws.Range("N84").Value2 = -38227.5
$ws.Range("H85").Value2 = 22916.666
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 22916.666
$ws.Range("K85").Value2 = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").Value2 = 68749.99800000001
$ws.Range("N85").Value2 = -71557.99800000001
$ws.Range("H87").Value2 = 27241.143
$ws.Range("I87").Value2 = 17850
$ws.Range("K87").Value2 = 53550
$ws.Range("M87").Value2 = -52302
$ws.Range("H88").Value2 = 9199
$ws.Range("J88").Value2 = 9199
$ws.Range("L88").Value2 = 27597
$ws.Range("N88").Value2 = -28453
$ws.Range("H90").Value2 = 27241.143
$ws.Range("I90").Value2 = 17850
$ws.Range("K90").Value2 = 160650
$ws.Range("M90").Value2 = -154410
$ws.Range("H91").Value2 = 9199
$ws.Range("J91").Value2 = 9199
$ws.Range("L91").Value2 = 27597
$ws.Range("N91").Value2 = -30561
$ws.Range("H95").Value2 = 17597.8
$ws.Range("I95").Value2 = 10000
$ws.Range("J95").Value2 = 19497.25
$ws.Range("K95").Value2 = 30000
$ws.Range("L95").Value2 = 58491.75
$ws.Range("M95").Value2 = -27941
$ws.Range("N95").Value2 = -62609.75
$ws.Range("H114").Value2 = 91696.37
$ws.Range("J114").Value2 = 334670
$ws.Range("L114").Value2 = 1004010
$ws.Range("N114").Value2 = -1010518
$ws.Range("H117").Value2 = 4233.4736
$ws.Range("J117").Value2 = 3504.1428
$ws.Range("L117").Value2 = 10512.4284
$ws.Range("N117").Value2 = -17396.4284
$ws.Range("H129").Value2 = 12821554
$ws.Range("I129").Value2 = 860.1111
$ws.Range("J129").Value2 = 41668110
$ws.Range("K129").Value2 = 2580.3333
$ws.Range("L129").Value2 = 125004330
$ws.Range("M129").Value2 = 2419.6667
$ws.Range("N129").Value2 = -125014330
$ws.Range("H132").Value2 = 3213
$ws.Range("I132").Value2 = 1000
$ws.Range("J132").Value2 = 3371.0715
$ws.Range("K132").Value2 = 9000
$ws.Range("L132").Value2 = 30339.6435
$ws.Range("M132").Value2 = -6470
$ws.Range("N132").Value2 = -35399.6435
$ws.Range("H133").Value2 = 6418.7144
$ws.Range("I133").Value2 = 6418.7144
$ws.Range("J133").Value2 = 0
$ws.Range("K133").Value2 = 19256.1432
$ws.Range("L133").Value2 = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value2 = -14196.1432
$ws.Range("H134").Value2 = 9393.888999999999
$ws.Range("I134").Value2 = 9672.666999999999
$ws.Range("J134").Value2 = 8000
$ws.Range("K134").Value2 = 29018.001
$ws.Range("L134").Value2 = 24000
$ws.Range("M134").Value2 = -23948.001
$ws.Range("N134").Value2 = -34140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value2 = 77.5
$ws.Range("I13").Value2 = 77.5
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 77.5
$ws.Range("L13").Value2 = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value2 = 61.5
$ws.Range("H70").Value2 = 7841
$ws.Range("I70").Value2 = 5376.5
$ws.Range("K70").Value2 = 5376.5
$ws.Range("M70").Value2 = -5106.5
$ws.Range("H73").Value2 = 7841
$ws.Range("I73").Value2 = 5376.5
$ws.Range("K73").Value2 = 5376.5
$ws.Range("M73").Value2 = -4440.5
$ws.Range("H80").Value2 = 6675.7144
$ws.Range("J80").Value2 = 10582.333
$ws.Range("L80").Value2 = 10582.333
$ws.Range("N80").Value2 = -12578.333
$ws.Range("H83").Value2 = 6675.7144
$ws.Range("J83").Value2 = 10582.333
$ws.Range("L83").Value2 = 52911.665
$ws.Range("N83").Value2 = -62895.665
$ws.Range("H97").Value2 = 1795.2
$ws.Range("I97").Value2 = 1731.3846
$ws.Range("K97").Value2 = 1731.3846
$ws.Range("M97").Value2 = -1235.3846
$ws.Range("H98").Value2 = 9820.75
$ws.Range("J98").Value2 = 9820.75
$ws.Range("L98").Value2 = 9820.75
$ws.Range("N98").Value2 = -15810.75
$ws.Range("H102").Value2 = 4134.846
$ws.Range("I102").Value2 = 2143.5217
$ws.Range("K102").Value2 = 2143.5217
$ws.Range("M102").Value2 = -521.5216999999998
$ws.Range("H113").Value2 = 6775.724
$ws.Range("I113").Value2 = 4929.048
$ws.Range("K113").Value2 = 4929.048
$ws.Range("M113").Value2 = -2759.048
$ws.Range("H122").Value2 = 6942.091
$ws.Range("I122").Value2 = 6026.778
$ws.Range("J122").Value2 = 11061
$ws.Range("K122").Value2 = 18080.334
$ws.Range("L122").Value2 = 33183
$ws.Range("M122").Value2 = -15630.334
$ws.Range("N122").Value2 = -38083
$ws.Range("H132").Value2 = 166670260
$ws.Range("I132").Value2 = 200002800
$ws.Range("J132").Value2 = 7500
$ws.Range("K132").Value2 = 600008400
$ws.Range("L132").Value2 = 22500
$ws.Range("M132").Value2 = -600005870
$ws.Range("N132").Value2 = -27560
$ws.Range("H134").Value2 = 38953.668
$ws.Range("J134").Value2 = 38953.668
$ws.Range("L134").Value2 = 116861.004
$ws.Range("N134").Value2 = -121931.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 7504.7856
$ws.Range("I7").Value2 = 5961.1577
$ws.Range("J7").Value2 = 10763.556
$ws.Range("K7").Value2 = 5961.1577
$ws.Range("L7").Value2 = 10763.556
$ws.Range("M7").Value2 = -5849.1577
$ws.Range("N7").Value2 = -10987.556
$ws.Range("H16").Value2 = 6383.0386
$ws.Range("I16").Value2 = 5728.1
$ws.Range("K16").Value2 = 5728.1
$ws.Range("M16").Value2 = -5558.1
$ws.Range("H22").Value2 = 0
$ws.Range("I22").Value2 = 0
$ws.Range("K22").Value2 = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value2 = 0
$ws.Range("I27").Value2 = 0
$ws.Range("K27").Value2 = 0
$ws.Range("M27").ClearContents()
$ws.Range("H43").Value2 = 2578749.8
$ws.Range("I43").Value2 = 149999
$ws.Range("J43").Value2 = 3388333.2
$ws.Range("K43").Value2 = 149999
$ws.Range("L43").Value2 = 3388333.2
$ws.Range("M43").Value2 = -149806
$ws.Range("N43").Value2 = -3388719.2
$ws.Range("H46").Value2 = 12821860
$ws.Range("I46").Value2 = 1071.45
$ws.Range("J46").Value2 = 26317426
$ws.Range("K46").Value2 = 1071.45
$ws.Range("L46").Value2 = 26317426
$ws.Range("M46").Value2 = -883.45
$ws.Range("N46").Value2 = -26317802
$ws.Range("H55").Value2 = 2269.4644
$ws.Range("I55").Value2 = 527.2353000000001
$ws.Range("J55").Value2 = 4962
$ws.Range("K55").Value2 = 527.2353000000001
$ws.Range("L55").Value2 = 4962
$ws.Range("M55").Value2 = -354.2353000000001
$ws.Range("N55").Value2 = -5308
$ws.Range("H82").Value2 = 3596.8
$ws.Range("I82").Value2 = 2244.4443
$ws.Range("J82").Value2 = 5625.3335
$ws.Range("K82").Value2 = 2244.4443
$ws.Range("L82").Value2 = 5625.3335
$ws.Range("M82").Value2 = -1883.4443
$ws.Range("N82").Value2 = -6347.3335
$ws.Range("H85").Value2 = 3596.8
$ws.Range("I85").Value2 = 2244.4443
$ws.Range("J85").Value2 = 5625.3335
$ws.Range("K85").Value2 = 2244.4443
$ws.Range("L85").Value2 = 5625.3335
$ws.Range("M85").Value2 = -996.4443000000001
$ws.Range("N85").Value2 = -8121.3335
$ws.Range("H87").Value2 = 69998
$ws.Range("J87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value2 = 69998
$ws.Range("J90").Value2 = 0
$ws.Range("L90").Value2 = 0
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value2 = 3414.0833
$ws.Range("I100").Value2 = 3369.9092
$ws.Range("J100").Value2 = 3900
$ws.Range("K100").Value2 = 3369.9092
$ws.Range("L100").Value2 = 3900
$ws.Range("M100").Value2 = -2828.9092
$ws.Range("N100").Value2 = -4982
$ws.Range("H122").Value2 = 5934.2173
$ws.Range("I122").Value2 = 5222.5
$ws.Range("J122").Value2 = 6313.8
$ws.Range("K122").Value2 = 15667.5
$ws.Range("L122").Value2 = 18941.4
$ws.Range("M122").Value2 = -13217.5
$ws.Range("N122").Value2 = -23841.4
$ws.Range("H126").Value2 = 7504.7856
$ws.Range("I126").Value2 = 5961.1577
$ws.Range("J126").Value2 = 10763.556
$ws.Range("K126").Value2 = 17883.4731
$ws.Range("L126").Value2 = 32290.668
$ws.Range("M126").Value2 = -15413.4731
$ws.Range("N126").Value2 = -37230.66800000001
$ws.Range("H132").Value2 = 5359.6113
$ws.Range("I132").Value2 = 5604.625
$ws.Range("J132").Value2 = 3399.5
$ws.Range("K132").Value2 = 16813.875
$ws.Range("L132").Value2 = 10198.5
$ws.Range("M132").Value2 = -14283.875
$ws.Range("N132").Value2 = -15258.5
$ws.Range("H135").Value2 = 154000
$ws.Range("J135").Value2 = 145000
$ws.Range("L135").Value2 = 145000
$ws.Range("N135").Value2 = -155140
$ws.Range("H136").Value2 = 26322886
$ws.Range("I136").Value2 = 41673612
$ws.Range("K136").Value2 = 125020836
$ws.Range("M136").Value2 = -125018286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value2 = 22656.428
$ws.Range("J26").Value2 = 23101.334
$ws.Range("L26").Value2 = 23101.334
$ws.Range("N26").Value2 = -23687.334
$ws.Range("H62").Value2 = 11166
$ws.Range("I62").Value2 = 9349.4
$ws.Range("J62").Value2 = 13436.75
$ws.Range("K62").Value2 = 9349.4
$ws.Range("L62").Value2 = 13436.75
$ws.Range("M62").Value2 = -8725.4
$ws.Range("N62").Value2 = -14684.75
$ws.Range("H65").Value2 = 11166
$ws.Range("I65").Value2 = 9349.4
$ws.Range("J65").Value2 = 13436.75
$ws.Range("K65").Value2 = 46747
$ws.Range("L65").Value2 = 67183.75
$ws.Range("M65").Value2 = -43627
$ws.Range("N65").Value2 = -73423.75
$ws.Range("H81").Value2 = 574.75
$ws.Range("I81").Value2 = 550
$ws.Range("J81").Value2 = 599.5
$ws.Range("K81").Value2 = 1100
$ws.Range("L81").Value2 = 1199
$ws.Range("M81").Value2 = -39
$ws.Range("N81").Value2 = -3321
$ws.Range("H84").Value2 = 574.75
$ws.Range("I84").Value2 = 550
$ws.Range("J84").Value2 = 599.5
$ws.Range("K84").Value2 = 5500
$ws.Range("L84").Value2 = 5995
$ws.Range("M84").Value2 = -196
$ws.Range("N84").Value2 = -16603
$ws.Range("H86").Value2 = 53438.332
$ws.Range("J86").Value2 = 53438.332
$ws.Range("L86").Value2 = 53438.332
$ws.Range("N86").Value2 = -55684.332
$ws.Range("H89").Value2 = 53438.332
$ws.Range("J89").Value2 = 53438.332
$ws.Range("L89").Value2 = 267191.66
$ws.Range("N89").Value2 = -278423.66
$ws.Range("H96").Value2 = 3216.5
$ws.Range("I96").Value2 = 2181.5454
$ws.Range("K96").Value2 = 2181.5454
$ws.Range("M96").Value2 = -808.5454
$ws.Range("H107").Value2 = 765.19354
$ws.Range("I107").Value2 = 470.85715
$ws.Range("K107").Value2 = 1412.57145
$ws.Range("M107").Value2 = 507.4285500000001
$ws.Range("H122").Value2 = 2735.4783
$ws.Range("I122").Value2 = 2548.353
$ws.Range("J122").Value2 = 3265.6667
$ws.Range("K122").Value2 = 7645.059
$ws.Range("L122").Value2 = 9797.000100000001
$ws.Range("M122").Value2 = -5195.059
$ws.Range("N122").Value2 = -14697.0001
$ws.Range("H126").Value2 = 4536.2
$ws.Range("I126").Value2 = 3561
$ws.Range("K126").Value2 = 10683
$ws.Range("M126").Value2 = -8213
$ws.Range("H132").Value2 = 14936.417
$ws.Range("I132").Value2 = 14535.429
$ws.Range("J132").Value2 = 15497.8
$ws.Range("K132").Value2 = 43606.287
$ws.Range("L132").Value2 = 46493.39999999999
$ws.Range("M132").Value2 = -41076.287
$ws.Range("N132").Value2 = -51553.39999999999
